$d = $word.ActiveDocument

$replacements = @(
    @{old = "946×8="; new = "517×2="},
    @{old = "845×4="; new = "708×9="},
    @{old = "706×9="; new = "620×9="},
    @{old = "955×5="; new = "827×9="},
    @{old = "279×8="; new = "772×2="},
    @{old = "253×8="; new = "218×9="},
    @{old = "367×4="; new = "421×4="},
    @{old = "323×5="; new = "358×7="},
    @{old = "935×8="; new = "187×9="},
    @{old = "801×9="; new = "956×5="},
    @{old = "936×5="; new = "421×3="},
    @{old = "692×3="; new = "573×7="},
    @{old = "892×7="; new = "745×9="},
    @{old = "241×8="; new = "694×7="},
    @{old = "977×4="; new = "318×5="},
    @{old = "446×8="; new = "388×5="},
    @{old = "217×2="; new = "307×6="},
    @{old = "332×3="; new = "132×7="},
    @{old = "370×8="; new = "367×6="},
    @{old = "376×4="; new = "455×9="},
    @{old = "422×9="; new = "303×9="},
    @{old = "931×6="; new = "886×9="},
    @{old = "150×9="; new = "593×5="},
    @{old = "493×3="; new = "890×4="},
    @{old = "839×7="; new = "688×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
